$d = $word.ActiveDocument

# Move to the very end of the document and add a new paragraph there,
# inheriting the paragraph/run formatting (list style, numbering, lang)
# from the paragraph that precedes it.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# Fill in the text of the newly created (last) paragraph.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rsquo = [char]0x2019
$newText = "Right now you" + $rsquo + "re working on this folder and later on you" + $rsquo + "ll update everything on the other folder. "
$newPara.Range.InsertBefore($newText)
